# Update "想去人数" (interest count) figures for several 北京-漫展信息 events.
# These numbers were refreshed when the static site was regenerated
# (gh-pages output at commit 456a3b4), so only column F values change.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsShow       = $wb.Worksheets.Item("演出")
$wsAll        = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$wsExhibition.Range("F3").Value  = 2851
$wsExhibition.Range("F6").Value  = 2479
$wsExhibition.Range("F11").Value = 45
$wsExhibition.Range("F13").Value = 7121
$wsExhibition.Range("F20").Value = 8360
$wsExhibition.Range("F34").Value = 2599
$wsExhibition.Range("F40").Value = 675
$wsExhibition.Range("F46").Value = 171
$wsExhibition.Range("F47").Value = 16

# 演出 (sheet2)
$wsShow.Range("F11").Value = 17

# 全部类型 (sheet4)
$wsAll.Range("F3").Value  = 2851
$wsAll.Range("F7").Value  = 2479
$wsAll.Range("F13").Value = 45
$wsAll.Range("F17").Value = 7121
$wsAll.Range("F23").Value = 8360
$wsAll.Range("F38").Value = 2599
$wsAll.Range("F43").Value = 675
$wsAll.Range("F49").Value = 171
